$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6564859299931365
$ws.Range("C2").Value = 0.721841688286301
$ws.Range("D2").Value = 0.6564859299931365
$ws.Range("E2").Value = 0.6558326373868439
$ws.Range("G2").Value = 0.7156651405486846
$ws.Range("I2").Value = 0.6753517827618257
$ws.Range("J2").Value = 0.6780141843971631
$ws.Range("K2").Value = 0.7261563106742928
$ws.Range("L2").Value = 0.6780141843971631
$ws.Range("M2").Value = 0.6852242934396159
$ws.Range("N2").Value = 0.7917867764813545
$ws.Range("O2").Value = 0.8049000220594653
$ws.Range("P2").Value = 0.7917867764813545
$ws.Range("Q2").Value = 0.7895324802895383
$ws.Range("S2").Value = 0.7781237330775052
$ws.Range("U2").Value = 0.7712840061605618
$ws.Range("V2").Value = 0.7768245252802564
$ws.Range("W2").Value = 0.7892638411582256
$ws.Range("X2").Value = 0.7768245252802564
$ws.Range("Y2").Value = 0.7754339125918925
$ws.Range("B3").Value = 0.8176847403340197
$ws.Range("C3").Value = 0.8306183837169406
$ws.Range("D3").Value = 0.8176847403340197
$ws.Range("E3").Value = 0.8185404329350711
$ws.Range("F3").Value = 0.8219400594829558
$ws.Range("G3").Value = 0.8349241610117154
$ws.Range("H3").Value = 0.8219400594829558
$ws.Range("I3").Value = 0.8232908288756373
$ws.Range("J3").Value = 0.8154655685197895
$ws.Range("K3").Value = 0.8255433940680035
$ws.Range("L3").Value = 0.8154655685197895
$ws.Range("M3").Value = 0.815189411956425
$ws.Range("N3").Value = 0.8412720201326926
$ws.Range("O3").Value = 0.8492533939538267
$ws.Range("P3").Value = 0.8412720201326926
$ws.Range("Q3").Value = 0.8410460689111131
$ws.Range("R3").Value = 0.8498284145504462
$ws.Range("S3").Value = 0.8571340348811713
$ws.Range("T3").Value = 0.8498284145504462
$ws.Range("U3").Value = 0.8494769476866167
$ws.Range("V3").Value = 0.8498055364905056
$ws.Range("W3").Value = 0.8570366339939139
$ws.Range("X3").Value = 0.8498055364905056
$ws.Range("Y3").Value = 0.8491655765090369
$ws.Range("B4").Value = 0.8411805078929306
$ws.Range("C4").Value = 0.8478428050927895
$ws.Range("D4").Value = 0.8411805078929306
$ws.Range("E4").Value = 0.84132206772719
$ws.Range("G4").Value = 0.8730950961824604
$ws.Range("I4").Value = 0.86723448644011
$ws.Range("J4").Value = 0.8412033859528713
$ws.Range("K4").Value = 0.8490249371817642
$ws.Range("L4").Value = 0.8412033859528713
$ws.Range("M4").Value = 0.8418780221405818
$ws.Range("N4").Value = 0.8520247083047356
$ws.Range("O4").Value = 0.8582545534475244
$ws.Range("P4").Value = 0.8520247083047356
$ws.Range("Q4").Value = 0.8516567262860322
$ws.Range("V4").Value = 0.8584534431480211
$ws.Range("W4").Value = 0.8643839902026158
$ws.Range("X4").Value = 0.8584534431480211
$ws.Range("Y4").Value = 0.8580478323653455
$ws.Range("B5").Value = 0.8369709448638755
$ws.Range("C5").Value = 0.8450048095849809
$ws.Range("D5").Value = 0.8369709448638755
$ws.Range("E5").Value = 0.834981754921556
$ws.Range("F5").Value = 0.8498284145504462
$ws.Range("G5").Value = 0.8593121634240288
$ws.Range("H5").Value = 0.8498284145504462
$ws.Range("I5").Value = 0.8482817846282467
$ws.Range("N5").Value = 0.8412262640128118
$ws.Range("O5").Value = 0.8529539208459198
$ws.Range("P5").Value = 0.8412262640128118
$ws.Range("Q5").Value = 0.8394540661465577
$ws.Range("R5").Value = 0.8455273392816289
$ws.Range("S5").Value = 0.85654401938693
$ws.Range("T5").Value = 0.8455273392816289
$ws.Range("U5").Value = 0.8435770137537075
$ws.Range("B6").Value = 0.8519560741249143
$ws.Range("C6").Value = 0.8588355878699682
$ws.Range("D6").Value = 0.8519560741249143
$ws.Range("E6").Value = 0.8514971511529369
$ws.Range("F6").Value = 0.8669869595058339
$ws.Range("G6").Value = 0.8717134422542763
$ws.Range("H6").Value = 0.8669869595058339
$ws.Range("I6").Value = 0.8669612660266498
$ws.Range("J6").Value = 0.8326698695950583
$ws.Range("K6").Value = 0.8461733539050206
$ws.Range("L6").Value = 0.8326698695950583
$ws.Range("M6").Value = 0.834593199735816
$ws.Range("N6").Value = 0.8583390528483186
$ws.Range("O6").Value = 0.8636111805228868
$ws.Range("P6").Value = 0.8583390528483186
$ws.Range("Q6").Value = 0.8577175509731683
$ws.Range("V6").Value = 0.8561885152139099
$ws.Range("W6").Value = 0.8621835881176961
$ws.Range("X6").Value = 0.8561885152139099
$ws.Range("Y6").Value = 0.8561439398225683
